# Apply fix: FHIR IG terminology and profile corrections
# - Set the "Experimental" metadata value (previously blank) to "false"
# - Update the "Date" metadata value to the new timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = Experimental property; its Value cell (B7) was empty, now "false"
# Force text storage (not boolean) - leading apostrophe marks it as literal text
$b7 = $ws.Range("B7")
$b7.Value = "'false"

# Re-apply the original (non-quote-prefixed) number format so the cell's
# style matches the rest of the body rows, restoring it to plain text style
$ws.Range("B6").Copy() | Out-Null
$b7.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 8 = Date property; update its Value cell (B8) to the new date/time
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
